$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 1.00000015925616
$ws.Range("G4").Value = 0.6107773305848241
$ws.Range("H4").Value = 1.63725814495874
$ws.Range("I4").Value = 0.002327894559130073
$ws.Range("J4").Value = 0.002149965708667878
$ws.Range("K4").Value = 0.3085486581549048
$ws.Range("L4").Value = 0.03071292489767075
$ws.Range("M4").Value = 0.0006180771888466552
$ws.Range("N4").Value = 0.2859652930637822
$ws.Range("O4").Value = 0.0001950151054188609
$ws.Range("P4").Value = 0.09720575390383601
